# Apply the edits described by the commit:
# "Added toggle to show passwords on manage students and manage instructors
#  page. Still deciding whether to keep manage advisors as requirements have
#  advisors many to many with departments."
#
# Concretely, in this workbook (backend/temporaryData.xlsx) the edit touches:
#   - the active/selected sheet & cell on a few sheets (UI state left over
#     from the author poking around in Excel while testing the app)
#   - the zoom level of the system_logs sheet
#   - the scroll position of the taken_data sheet
#   - a batch of "grade" values on taken_data that got rounded from one
#     decimal place down to whole numbers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# taken_data: round several grade values (column C) to whole numbers
# ---------------------------------------------------------------------
$wsTaken = $wb.Worksheets.Item("taken_data")

$wsTaken.Range("C3").Value  = 3
$wsTaken.Range("C4").Value  = 3
$wsTaken.Range("C6").Value  = 3
$wsTaken.Range("C7").Value  = 3
$wsTaken.Range("C8").Value  = 3
$wsTaken.Range("C10").Value = 3
$wsTaken.Range("C12").Value = 2
$wsTaken.Range("C15").Value = 2
$wsTaken.Range("C16").Value = 2
$wsTaken.Range("C19").Value = 3
$wsTaken.Range("C23").Value = 3
$wsTaken.Range("C27").Value = 3
$wsTaken.Range("C28").Value = 2
$wsTaken.Range("C33").Value = 3
$wsTaken.Range("C35").Value = 3

# taken_data view: scroll down a bit (topLeftCell A7 -> A13) while keeping
# the previous selection at H31
$wsTaken.Activate()
$wsTaken.Range("H31").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# system_logs: zoom out from 160% to 130%
# ---------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("system_logs")
$wsLogs.Activate()
$wsLogs.Range("F13").Select()
$excel.ActiveWindow.Zoom = 130

# ---------------------------------------------------------------------
# major_data: no longer the tab that is active/selected when reopening
# (tabSelected flag moves to student_data below); keep its own selection
# as-is (I18)
# ---------------------------------------------------------------------
$wsMajor = $wb.Worksheets.Item("major_data")
$wsMajor.Range("I18").Select()

# ---------------------------------------------------------------------
# student_data: becomes the active sheet/tab, with a new active cell/
# selection at J28
# ---------------------------------------------------------------------
$wsStudent = $wb.Worksheets.Item("student_data")
$wsStudent.Activate()
$wsStudent.Range("J28").Select()
